# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Sheet "Rules", cell B11 changes from the text "R40" to the text "1".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
# Force the numeric-looking literal to be stored as text (leading apostrophe),
# matching the original value which was also stored as a shared string.
$cell.Value = "'1"
